$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on price cells whose new values would
# otherwise be auto-interpreted as numbers (losing the original
# text/string cell type used throughout column D).
$textCells = 'D5,D11,D12,D14,D17,D22,D23,D25,D26,D28,D30,D31,D32,D33,D34,D37,D38,D39,D49,D51'.Split(',')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range('D2').Value = '64.148.33'
$ws.Range('E2').Value = '  -0.89%  '

# Row 3
$ws.Range('D3').Value = '3.324.48'
$ws.Range('E3').Value = '  -1.64%  '

# Row 4
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$ws.Range('D5').Value = '551.43'
$ws.Range('E5').Value = '  -1.16%  '

# Row 6
$ws.Range('E6').Value = '  -2.36%  '

# Row 7
$ws.Range('E7').Value = '  +1.09%  '

# Row 8
$ws.Range('E8').Value = '  +0.08%  '

# Row 9
$ws.Range('D9').Value = '3.314.00'
$ws.Range('E9').Value = '  -1.72%  '

# Row 10
$ws.Range('E10').Value = '  +5.46%  '

# Row 11
$ws.Range('D11').Value = '0.636'
$ws.Range('E11').Value = '  +0.89%  '

# Row 12
$ws.Range('D12').Value = '53.20'
$ws.Range('E12').Value = '  -3.22%  '

# Row 13
$ws.Range('E13').Value = '  +1.50%  '

# Row 14
$ws.Range('D14').Value = '9.05'
$ws.Range('E14').Value = '  -0.56%  '

# Row 15
$ws.Range('D15').Value = '3.854.44'
$ws.Range('E15').Value = '  -1.08%  '

# Row 16
$ws.Range('E16').Value = '  +1.91%  '

# Row 17
$ws.Range('D17').Value = '18.06'
$ws.Range('E17').Value = '  -2.04%  '

# Row 18
$ws.Range('D18').Value = '3.335.66'
$ws.Range('E18').Value = '  -0.66%  '

# Row 19
$ws.Range('D19').Value = '63.940.54'
$ws.Range('E19').Value = '  -1.05%  '

# Row 20
$ws.Range('E20').Value = '  -1.86%  '

# Row 21
$ws.Range('E21').Value = '  -0.38%  '

# Row 22
$ws.Range('D22').Value = '446.52'
$ws.Range('E22').Value = '  +3.23%  '

# Row 23
$ws.Range('D23').Value = '4.95'
$ws.Range('E23').Value = '  +0.84%  '

# Row 24
$ws.Range('E24').Value = '  -2.30%  '

# Row 25
$ws.Range('D25').Value = '86.79'
$ws.Range('E25').Value = '  +2.87%  '

# Row 26
$ws.Range('D26').Value = '13.79'
$ws.Range('E26').Value = '  +4.16%  '

# Row 27
$ws.Range('E27').Value = '  +0.51%  '

# Row 28
$ws.Range('D28').Value = '10.59'
$ws.Range('E28').Value = '  -1.97%  '

# Row 29
$ws.Range('E29').Value = '  -2.57%  '

# Row 30
$ws.Range('D30').Value = '30.76'
$ws.Range('E30').Value = '  +2.95%  '

# Row 31
$ws.Range('D31').Value = '6.49'
$ws.Range('E31').Value = '  -2.94%  '

# Row 32
$ws.Range('D32').Value = '62.25'
$ws.Range('E32').Value = '  +6.45%  '

# Row 33
$ws.Range('D33').Value = '11.35'
$ws.Range('E33').Value = '  -1.41%  '

# Row 34
$ws.Range('D34').Value = '570.35'
$ws.Range('E34').Value = '  -0.34%  '

# Row 35
$ws.Range('E35').Value = '  -2.06%  '

# Row 36
$ws.Range('E36').Value = '  -0.06%  '

# Row 37
$ws.Range('B37').Value = 'Stacks'
$ws.Range('C37').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D37').Value = '3.56'
$ws.Range('E37').Value = '  +1.23%  '

# Row 38
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').Value = '0.141'
$ws.Range('E38').Value = '  -1.40%  '

# Row 39
$ws.Range('D39').Value = '35.12'
$ws.Range('E39').Value = '  -2.19%  '

# Row 40
$ws.Range('E40').Value = '  -1.38%  '

# Row 41
$ws.Range('D41').Value = '0.0₃0727'
$ws.Range('E41').Value = '  -4.51%  '

# Row 42
$ws.Range('D42').Value = '3.057.44'
$ws.Range('E42').Value = '  -2.05%  '

# Row 43
$ws.Range('E43').Value = '  +0.19%  '

# Row 44
$ws.Range('E44').Value = '  -4.36%  '

# Row 45
$ws.Range('E45').Value = '  -2.96%  '

# Row 46
$ws.Range('E46').Value = '  +2.28%  '

# Row 47
$ws.Range('E47').Value = '  -1.68%  '

# Row 48
$ws.Range('E48').Value = '  +0.08%  '

# Row 49
$ws.Range('D49').Value = '142.25'
$ws.Range('E49').Value = '  +5.39%  '

# Row 50
$ws.Range('E50').Value = '  -2.88%  '

# Row 51
$ws.Range('D51').Value = '8.15'
$ws.Range('E51').Value = '  -1.68%  '
